$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (row 2, column B) no longer returns the Cohort
# column - drop the trailing `coalesce(co.cohort_description, '') AS
# `Cohort`` clause from the Cypher query text (B3/B4 - the Samples/Files
# tab queries - are untouched).
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
    WHERE f.file_format IN ["doc"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# Matching layout refresh that Excel performed when it rewrapped the
# shortened query text: the row heights for the three query rows were
# recomputed.
$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 259.2

# Selection moved back up to B2 (and the view scrolled back so row 1 is
# visible, instead of being scrolled down to row 4).
$ws.Range("B2").Select() | Out-Null
